$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Dodanie podziału treningu na części" — split the training session into
# parts ("Duża Gra" / "Mała Gra"). A new "Trening" column (F) is added, the
# existing rows are re-labelled "Mała Gra" and pushed down, and a brand new
# block of "Duża Gra" rows is inserted above them (plus one extra "Mała Gra"
# row at the end).
# ---------------------------------------------------------------------------

# New header for column F — same bold/centered/bordered style as the rest
# of the header row (style index 1, reused from E1 via a format-only copy).
$ws.Cells.Item(1, 6).Value = "Trening"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Full final dataset for rows 2-13:
#   Timestamp(serial), Seconds, Velocity, Acceleration_SMA, Velocity_Bin, Trening
$data = @(
    @(45684.59115486111, 475.7,  10.76, 3.104219096047534, "10-15", "Duża Gra"),
    @(45684.59195231481, 544.6,  14.54, 2.930859497615269, "10-15", "Duża Gra"),
    @(45684.59379837963, 704.1,  11.77, 2.81770658493042,  "10-15", "Duża Gra"),
    @(45684.59115254629, 475.5,  9.720000000000001, 3.056578295571462, "5-10", "Duża Gra"),
    @(45684.59379606482, 703.9,  9.26,  2.633463025093079, "5-10",  "Duża Gra"),
    @(45684.59401944444, 723.2,  9.619999999999999, 2.420666234833852, "5-10", "Duża Gra"),
    @(45684.59746736111, 1021.1, 13.91, 4.529195376804894, "10-15", "Mała Gra"),
    @(45684.59890601852, 1145.4, 13.96, 4.76598743030003,  "10-15", "Mała Gra"),
    @(45684.60191990741, 1405.8, 14.85, 4.633814913885933, "10-15", "Mała Gra"),
    @(45684.59746388889, 1020.8, 8.66,  3.780123114585876, "5-10",  "Mała Gra"),
    @(45684.60191643518, 1405.5, 9.34,  3.636793834822519, "5-10",  "Mała Gra"),
    @(45684.60288518519, 1489.2, 8.1,   3.464750153677804, "5-10",  "Mała Gra")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Apply the datetime number format to the Timestamp column for every data
# row (the first assignment, lowercase, then the real uppercase one —
# mirrors how the format ended up registered twice in the workbook).
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2:A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
